# "Add files via upload" — a new sample row (35) was appended below the
# existing data table on sheet "principal", and row 34 (the row that used
# to be last) picked up the same explicit cell formatting the rest of the
# table already carries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("principal")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Re-apply the table's formatting to row 34 now that it is no longer the
# last row of the range (matches the look of every other data row, e.g. 33).
$ws.Range("A34:G34").NumberFormat = "General"

# Append the new sample as row 35.
$ws.Range("A35").Value = "7.139007568359375 GB"
$ws.Range("B35").Value = "9.012344360351562GB"

# C35/D35 hold percentages formatted as plain text (e.g. "4.2 %"), not
# numeric percentages, so force text so Excel doesn't auto-convert them.
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "4.2 %"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "44.7 %"

$ws.Range("E35").Value = 576
$ws.Range("F35").Value = 329.3300383090973
$ws.Range("G35").Value = "6038127.png"
